# "Especificacao do caso de uso Tela Inicial"
# Adds requirement #14 / use case #13 ("Exibir Tela Inicial") as a new
# row at the bottom of the Requisitos e Casos de Uso table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table row (row 18): REQ# / Descricao / # / CASO DE USO
$ws.Range("A18").Value = 14
$ws.Range("B18").Value = "Exibir os últimos/próximos eventos na tela inicial"
$ws.Range("C18").Value = 13
$ws.Range("D18").Value = "Exibir Tela Inicial"

# Restore the scrolled view / active-cell selection left by the edit.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D19").Select()
